$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 67
$ws.Range("I2").Value = 164
$ws.Range("J2").Value = 684
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 196
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 132
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 92
$ws.Range("T2").Value = 107
$ws.Range("V2").Value = 1132
$ws.Range("X2").Value = 1097
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 6
